$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows with new lab member data
$ws.Range("A2").Value = "sasaki"
$ws.Range("B2").Value = "佐々木一郎"
$ws.Range("A3").Value = "sato"
$ws.Range("B3").Value = "佐藤二郎"
$ws.Range("A4").Value = "tanaka"
$ws.Range("B4").Value = "田中三郎"

# Remove the now-unused rows 5 and 6
$ws.Range("A5:B6").Delete()

# Update selection to match the new active cell
$ws.Range("B4").Select()
